$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New admission-list rows appended after the previous last row (224), for the
# 13:11 21.07.2025 batch of applicants. Columns: A F.I.Sh, B Yo'nalish,
# C Ta'lim tili, D Ta'lim shakli, E Passport, F JSHIR, G Viloyat, H Tuman/shahar,
# I/J phone numbers, K application date.
$rows = @(
    ,@(225, "Baxtiyorov ulugbek", "Yurisprudensiya", "Rus tili", "Kunduzgi", "AD5227398", "50811075740064", "Toshkent shahri", "Yakkasaroy tumani", "998932407777", "+998908771777", "2025-07-18")
    ,@(226, "Ibragimova Zulxumor Atabek Qizi", "Yurisprudensiya", "O'zbek tili", "Kunduzgi", "AD9594019", "61312067130051", "Xorazm viloyati", "Urganch shahri", "998981118202", "+998957115151", "2025-07-18")
    ,@(227, "Ibragimova Dilafruz Atabek qizi", "Yurisprudensiya", "O'zbek tili", "Kunduzgi", "AD9535989", "60912077130051", "Xorazm viloyati", "Urganch tumani", "998910141284", "+998946461284", "2025-07-18")
    ,@(228, "Qosimov Hikmatilla Lutfulla ogli", "Yurisprudensiya", "O'zbek tili", "Kunduzgi", "AD1296541", "30705910211304", "Toshkent shahri", "Yashnaobod tumani", "998974573333", "+998974573333", "2025-07-19")
    ,@(229, "Abduxoliqov Iskandarbek Davronjon o'gli", "Yurisprudensiya", "O'zbek tili", "Kunduzgi", "AD4537927", "52905076920012", "Fargona viloyati", "Buvayda tumani", "998908561835", "+998975301040", "2025-07-19")
    ,@(230, "Masharifov masharif murod ogli", "Yurisprudensiya", "O'zbek tili", "Kunduzgi", "AD4181861", "51212057130018", "Xorazm viloyati", "Urganch tumani", "998200056838", "+998918686838", "2025-07-19")
    ,@(231, "Xoliqov Nuriddin Damirovich", "Mehnat muhofazasi va texnika xavfsizligi", "O'zbek tili", "Kunduzgi", "AD9119669", "51710085360046", "Buxoro viloyati", "Shofirkon tumani", "998509008511", "+998501503580", "2025-07-20")
    ,@(232, "Hasansher Norboboyev Jaxongir o'g'li", "Yurisprudensiya", "O'zbek tili", "Kunduzgi", "AD0949619", "50307055680025", "Qashqadaryo viloyati", "Kitob tumani", "998770148278", "+998770148278", "2025-07-21")
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K")

foreach ($row in $rows) {
    $r = $row[0]

    # F (JSHIR), I/J (phone numbers) and K (date) all look numeric/date-like;
    # force text format first so Excel stores them as literal text, matching
    # the source data, instead of silently converting to numbers/dates.
    $ws.Range("F" + $r).NumberFormat = "@"
    $ws.Range("I" + $r + ":J" + $r).NumberFormat = "@"
    $ws.Range("K" + $r).NumberFormat = "@"

    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $row[$i + 1]
    }
}
